# Split "DD-DD | rest of heading" Heading3 runs into:
#   1) "DD-DD"  (bold)
#   2) " "      (plain)
#   3) "rest of heading" (plain)
# i.e. drop the " | " separator in favor of a single space, and bold the
# leading date code.

$d = $word.ActiveDocument

# pkg:package XML template used to splice three runs (with explicit rPr)
# into a paragraph range via Range.InsertXML, so formatting (w:b + w:bCs)
# round-trips exactly as authored.
$xmlTemplate = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">{0}</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">{1}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Collect target paragraph indices first (mutating while iterating the
# live collection can shift indices/ranges).
$targets = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 3") {
        $t = $p.Range.Text
        if ($t -match "^(\d\d-\d\d) \| (.+?)[\r\x07]?$") {
            [void]$targets.Add($i)
        }
    }
}

# Apply edits back-to-front so earlier paragraph ranges stay valid as we
# shrink/replace each target paragraph's content.
for ($j = $targets.Count - 1; $j -ge 0; $j--) {
    $i = $targets[$j]
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range
    $t = $full.Text
    if ($t -match "^(\d\d-\d\d) \| (.+?)[\r\x07]?$") {
        $datePart = $matches[1]
        $restPart = $matches[2]
        $start = $full.Start
        $end = $full.End
        # exclude the trailing paragraph mark
        $r = $d.Range($start, $end - 1)
        $frag = $xmlTemplate -f $datePart, $restPart
        $r.InsertXML($frag)
    }
}

Write-Output "Updated $($targets.Count) heading(s)"
